# The commit completes the calc_reactions / calc_pier_reactions tests for
# span 2 by changing the axle-location input `x` (cell C2) from 75 to 125.
# All of the dependent `abs axle location` formulas in column D
# (IF($C$3="ltr",$C$2-C5, $C$2+C5), etc.) recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 125

# Cosmetic view-state touch-ups from the original edit (scroll the active
# window back to the top-left cell and restore the default tab-bar ratio).
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.TabRatio = 0.6

$wb.Application.CalculateFull()
